$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 20:52"

# Swap the Asturias / Gipuzkoa rows (row 23 and row 24) and refresh their data
$ws.Range("A23").Value = "Asturias"
$ws.Range("B23").Value = 1958
$ws.Range("C23").Value = 434
$ws.Range("D23").Value = 1375
$ws.Range("E23").Value = 149

$ws.Range("A24").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("B24").Value = 1955
$ws.Range("C24").Value = 4867
$ws.Range("D24").Value = 5101
$ws.Range("E24").Value = 130
